$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.594.99"
$ws.Range('D2').Style = "Normal"
$ws.Range('D3').Value = "'1.647.40"
$ws.Range('D3').Style = "Normal"
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'212.67"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('E6').Value = '  +4.79%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = "'23.58"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -1.50%  '
$ws.Range('E9').Value = '  -1.63%  '
$ws.Range('E10').Value = '  -1.15%  '
$ws.Range('E11').Value = '  +1.13%  '
$ws.Range('D12').Value = "'1.880.88"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('D13').Value = "'1.649.50"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.36%  '
$ws.Range('D14').Value = "'0.586"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +4.24%  '
$ws.Range('E15').Value = '  -2.34%  '
$ws.Range('D16').Value = "'64.46"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.02%  '
$ws.Range('D17').Value = "'27.561.69"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = "'231.99"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -3.57%  '
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('D20').Value = "'7.57"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E22').Value = '  -2.97%  '
$ws.Range('D23').Value = "'9.74"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.67%  '
$ws.Range('D25').Value = "'149.03"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.31%  '
$ws.Range('E26').Value = '  -2.53%  '
$ws.Range('E27').Value = '  +1.48%  '
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('E29').Value = '  -3.90%  '
$ws.Range('E30').Value = '  -2.24%  '
$ws.Range('E31').Value = '  -2.95%  '
$ws.Range('E32').Value = '  -0.44%  '
$ws.Range('D33').Value = "'3.18"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.48%  '
$ws.Range('D34').Value = "'1.426.53"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('E35').Value = '  +2.86%  '
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').Value = "'0.571"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('E38').Value = '  -4.01%  '
$ws.Range('E39').Value = '  -2.13%  '
$ws.Range('E40').Value = '  -1.20%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('E42').Value = '  +2.32%  '
$ws.Range('E43').Value = '  +2.40%  '
$ws.Range('E44').Value = '  +1.20%  '
$ws.Range('D45').Value = "'65.16"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.61%  '
$ws.Range('D46').Value = "'1.790.18"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('E47').Value = '  -1.20%  '
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('E49').Value = '  +1.03%  '
$ws.Range('D50').Value = "'0.0997"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.31%  '
$ws.Range('D51').Value = "'7.77"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.99%  '
